$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Create new row 31 by cloning the formatting of row 3 (same column style pattern:
#    A=ALL_PAGES, B=END_TO_END, C=style11(+red font), D=wrapped description, E=RunMode, F=High)
$ws.Range("A3:F3").Copy($ws.Range("A31:F31"))

# 2. Fill in the new testcase content for row 31
$ws.Range("C31").Value = "TC40_Adding_MultipleItems_QuickOrder"
$ws.Range("D31").Value = "1. more than 10 Items should be added to cart`n2. Appropriate inventory message must be displayed"
$ws.Range("E31").Value = "Yes"

# 3. Highlight the new testcase id in red
$ws.Range("C31").Font.Color = 255

# 4. Match row height used by similar two-line wrapped rows
$ws.Rows("31").RowHeight = 30

# 5. Flip RunMode to "No" for all of the pre-existing testcases
$ws.Range("E2:E30").Value = "No"

# 6. Update the view to reflect the newly added row
$ws.Application.ActiveWindow.ScrollRow = 22
$ws.Range("D34").Select()
